$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Bai 27")

# --- Add the new rows of content to Sheet1 (rows 18-20, 22-31; row 21 left blank) ---
$ws1.Range("A18").Value = "Bài 28"
$ws1.Range("A19").Value = "Kết quả 517"
$ws1.Range("A20").Value = "Bài 32"

$ws1.Range("A22").Value = "Bài 33"
$ws1.Range("A23").Value = "Không hiểu bài 8 nào"
$ws1.Range("A24").Value = "bài 34"
$ws1.Range("A25").Value = "Không hiểu"
$ws1.Range("A26").Value = "bài 35 "
$ws1.Range("A27").Value = "Hoang Mang"
$ws1.Range("A28").Value = "Bài 36"
$ws1.Range("A29").Value = "Chưa hiểu"
$ws1.Range("A30").Value = "Bài 37 "
$ws1.Range("A31").Value = "ok"

# --- Switch the active sheet / selection from "Bai 27" to "Sheet1" ---
$ws1.Activate()
$ws1.Range("H27").Select()
